# The <id> paragraph ("<id>p145v_2</id>") was split across three runs:
#   "<id>p145v_"  (w:color=000000)
#   "2"           (no explicit color -> automatic)
#   "</id>"       (w:color=000000)
# The commit collapses these into a single run of "<id>p145v_2</id>"
# (keeping the surrounding black-color formatting). A literal Find/Replace
# over that exact text re-creates the match as one uniformly formatted run
# (taking on the leading run's formatting) while leaving the rest of the
# paragraph (including the trailing empty run) untouched.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p145v_2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p145v_2</id>", 2)
